$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (dates are Excel serial numbers for
# 2021-05-10 .. 2021-05-13, matching the existing column A date format)
$newRows = @(
    @{ Row = 252; A = 44326; B = 0; C = 6; D = 47.11425206124853 },
    @{ Row = 253; A = 44327; B = 1; C = 7; D = 54.96662740478995 },
    @{ Row = 254; A = 44328; B = 1; C = 8; D = 62.81900274833137 },
    @{ Row = 255; A = 44329; B = 1; C = 9; D = 70.67137809187278 }
)

foreach ($r in $newRows) {
    $rowIdx = $r.Row

    # Clone formatting of the row above (column A carries the date style)
    # by copy/paste-special of formats only, then fill in the values.
    $srcA = $ws.Cells.Item($rowIdx - 1, 1)
    $dstA = $ws.Cells.Item($rowIdx, 1)
    $srcA.Copy() | Out-Null
    $dstA.PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    $dstA.Value = $r.A
    $ws.Cells.Item($rowIdx, 2).Value = $r.B
    $ws.Cells.Item($rowIdx, 3).Value = $r.C
    $ws.Cells.Item($rowIdx, 4).Value = $r.D
}
